$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.290.16"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "1.855.17"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "'314.70"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").Value = "'0.4607"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "'0.3709"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").Value = "'0.07299"
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("D10").Value = "'0.8892"
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").Value = "'20.06"
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("D12").Value = "'0.07823"
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").Value = "1.819.09"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").Value = "'5.390"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").Value = "'6.522"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "'91.49"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "'0.000008926"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "27.304.01"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "'5.115"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").Value = "'10.55"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "2.062.33"
$ws.Range("E24").Value = "  +2.94%  "
$ws.Range("D25").Value = "'1.926"
$ws.Range("E25").Value = "  +4.37%  "
$ws.Range("D26").Value = "'152.11"
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("D27").Value = "'18.45"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").Value = "'2.059"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").Value = "'116.02"
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("D30").Value = "'5.070"
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("D31").Value = "'0.08827"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").Value = "'0.7739"
$ws.Range("E32").Value = "  +5.90%  "
$ws.Range("D33").Value = "'3.092"
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("D34").Value = "'1.173"
$ws.Range("E34").Value = "  +3.78%  "
$ws.Range("D35").Value = "'4.516"
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("D36").Value = "'2.747"
$ws.Range("E36").Value = "  +13.06%  "
$ws.Range("D37").Value = "'1.078"
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("D38").Value = "'0.01955"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").Value = "'0.05266"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").Value = "'2.955"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").Value = "'7.051"
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").Value = "'0.5128"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "'0.1640"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("D44").Value = "'8.418"
$ws.Range("E44").Value = "  +2.58%  "
$ws.Range("D45").Value = "'0.4792"
$ws.Range("D46").Value = "'10.33"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").Value = "'102.53"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").Value = "'1.643"
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("D50").Value = "'0.06219"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "'65.85"
$ws.Range("E51").Value = "  +2.00%  "
